$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark from the "Rockhopper penguin" paragraph ---
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# --- 2. Append the new animal-info content at the end of the document body ---
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document><w:body>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>PENGUIN</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Average</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> Weight: 30kg - 66lb</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Average Height: 115cm - 3.8ft</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Breeding Season: April - December - Emperors breed in the depth of the Antarctic winter. The average temperature is around -20°C (-4°F) falling as low as -50°C (- 58°F) and with winds that gust up to 200km per hour (124mph).</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Reproduction: Emperor Penguins breed almost exclusively on sea-ice near to the coast of Antarctica, rarely on land, the breeding colonies can be up to 200 km from the open sea.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Estimated world population: - 595,000, stable.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Feeding: Eat mainly fish in Antarctica, also krill and squid which can be more important in some places.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Conservation status: Near </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>threatened..</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Distributionn</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: Circumpolar, breed right up against the continent.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Predators: Leopard seals and killer whales- main predators of adult birds. Skuas - prey on eggs and chicks.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>SEAL</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Weight: Males 130 to 210 kg / females 22 - 55 kg.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Length: Males typically 1.8 m / females 1.2 - 1.4 m</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Breeding Season: Females arrive on breeding beaches a day or two before giving birth in November and December. Males are already in residence arriving a month or so earlier and have established territories. The small females stay with their pups for 6-7 days before returning to sea to feed, they come into oestrus at this time and mate again, they nurse the pups for about 4 months.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Estimated world population: - 5 to 7.5 million, population increasing.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Feeding: Fur seals will usually feed on krill where it is available and then on fish and squid when it is not. The availability varies with season and with location. They are also known to feed on penguins to a lesser extent.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Diving: Antarctic Fur seals are not great divers when compared to other seals, the longest dives are about 10 minutes, average dives are to 30-40m for around 2 minutes. Most dives are shallower than this when they do most of their feeding. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Conservation status: Least concern. Protected by the Antarctic Treaty and the Convention for the Conservation of Antarctic Seals.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Distribution: Circumpolar where there are seasonally ice-free islands, they are rarely found hauled out on ice preferring land. Mainly found just south or sometimes north of the Antarctic convergence, a sub-Antarctic species rather than Antarctic.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Predators: Leopard seals, especially of pups, in the South Shetland Islands predation levels are thought to be enough to cause a population decline.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>WHALE</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Whales are huge, but elusive and difficult to see which adds to their mystery and fascination. They are highly intelligent animals with an elaborate social life, no possessions and the complete freedom of movement in three dimensions. Is it any wonder that they are such popular and fascinating animals? - maybe we just want to be like them.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Whales belong to the group of mammals called Cetaceans, they are a part of this group along with dolphins and porpoises.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Whales are mammals as are humans, dogs, cats, elephants and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>anguantibos</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> amongst others. This means that they are not fish. They breathe air and so must return to the surface at regular intervals to get a breath. They give birth to live young that stay with the mother for over a year and feed on milk produced by the mother. </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>ARTIC FOX</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Average Weight: 3 to 8kg (6.5-17 lbs)</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Average Length: 75 to 110 cm long (2.3-3.5 feet) including a tail of around 30cm (12 inches), 20 to 30cm (9-12 inches) tall at the shoulder, females slightly smaller than males</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Predators: Polar bears, wolverines, red foxes and golden eagles. They have been and continue to be trapped for their thick winter coats in particular</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Estimated world population: - Several hundred thousand, wide fluctuations </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>as a result of</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> varying prey numbers and a rapid ability to reproduce.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Feeding: A wide range of foods, the main prey is lemmings, they will hunt and catch other small animals and will also scavenge food from beneath sea-bird colonies on cliffs and left </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>overs</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> from predators such as polar bears. They will take bird eggs where possible from tundra nesting birds, though are not entirely carnivorous eating berries and seaweed when available.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">A family of foxes can get through several dozen lemmings in a day. They will eat young ringed seals when they are vulnerable in the snow den shortly after they are born in the same manner that they attack lemmings beneath the snow, detecting them by sound and then jumping on and punching through the covering snow layer. </w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertXML($xml)
